$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the old example row (row 2); row 3 shifts up to become row 2 ---
$ws.Rows.Item(2).Delete()

# --- 2. Update header row (row 1) text with appended instructional lines ---
$ws.Cells.Item(1,1).Value = "구분`n신규 = 1`n경력 = 2"
$ws.Cells.Item(1,2).Value = "이름"
$ws.Cells.Item(1,3).Value = "주민등록번호`n000000-0000000"
$ws.Cells.Item(1,4).Value = "교육이수번호`n0000000000"
$ws.Cells.Item(1,5).Value = "경력시작일`nYYYY-MM-DD"
$ws.Cells.Item(1,6).Value = "경력종료일`nYYYY-MM-DD"

# --- 3. Format header row: wrap text, center, fill+border, taller row ---
$headerRange = $ws.Range("A1:F1")
$headerRange.WrapText = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4108    # xlCenter
$ws.Rows.Item(1).RowHeight = 52.2

# B1 ("이름") keeps no wrap (narrower content), restore
$ws.Cells.Item(1,2).WrapText = $false

# --- 4. Update data row (now row 2, former row 3) ---
$ws.Cells.Item(2,2).Value = "홍길동"

# --- 5. Column widths ---
$ws.Columns.Item(5).ColumnWidth = 15
$ws.Columns.Item(6).ColumnWidth = 15.296875

# --- 6. Selection ---
$ws.Range("F6").Select()
